$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.00001182803375754382
$ws.Range("D2").Value = 0.01030109339467344
$ws.Range("E2").Value = 0.049748111630187
$ws.Range("F2").Value = 1.08748851459336
$ws.Range("G2").Value = 1.050800046505458
$ws.Range("H2").Value = 0.8060506926459254
$ws.Range("I2").Value = 0.8258888374792406
$ws.Range("M2").Value = 1.225885884648903
$ws.Range("N2").Value = 1.62366883687983
$ws.Range("C3").Value = 0.00001681388636876058
$ws.Range("D3").Value = 0.009400443839112427
$ws.Range("E3").Value = 0.04988009129297666
$ws.Range("F3").Value = 1.0055234076803
$ws.Range("G3").Value = 0.9543758343926925
$ws.Range("H3").Value = 0.7676947923526711
$ws.Range("I3").Value = 0.7630507381648073
$ws.Range("M3").Value = 1.079393240439117
$ws.Range("N3").Value = 1.495714242357394
$ws.Range("C4").Value = 0.00002043077137448002
$ws.Range("D4").Value = 0.008860912165353341
$ws.Range("E4").Value = 0.0500219451821593
$ws.Range("F4").Value = 0.9560577475651399
$ws.Range("G4").Value = 0.8959646918607405
$ws.Range("H4").Value = 0.7448086268890677
$ws.Range("I4").Value = 0.7251253590535782
$ws.Range("M4").Value = 0.9894087817867216
$ws.Range("N4").Value = 1.41737298981576
$ws.Range("C5").Value = 0.00002203776956011261
$ws.Range("D5").Value = 0.008644319369224718
$ws.Range("E5").Value = 0.05009511842064285
$ws.Range("F5").Value = 0.9361124065218149
$ws.Range("G5").Value = 0.8723561596830507
$ws.Range("H5").Value = 0.7356471872757311
$ws.Range("I5").Value = 0.7098324735288202
$ws.Range("M5").Value = 0.9527298002085303
$ws.Range("N5").Value = 1.38550805751683
$ws.Range("C6").Value = 0.00002231245790906833
$ws.Range("D6").Value = 0.008608548460998122
$ws.Range("E6").Value = 0.0501081990536445
$ws.Range("F6").Value = 0.9328131979803089
$ws.Range("G6").Value = 0.8684475791007458
$ws.Range("H6").Value = 0.734135829369734
$ws.Range("I6").Value = 0.707302790760437
$ws.Range("M6").Value = 0.9466387003686521
$ws.Range("N6").Value = 1.380220626212946
$ws.Range("C7").Value = 0.00002045191408073244
$ws.Range("D7").Value = 0.008857978028878222
$ws.Range("E7").Value = 0.05002286970476
$ws.Range("F7").Value = 0.9557879036356525
$ws.Range("G7").Value = 0.8956455181891272
$ws.Range("H7").Value = 0.7446844079383936
$ws.Range("I7").Value = 0.7249184622074125
$ws.Range("M7").Value = 0.9889141552222043
$ws.Range("N7").Value = 1.416943001403212
$ws.Range("C8").Value = 0.00001342783489088895
$ws.Range("D8").Value = 0.009987682766634975
$ws.Range("E8").Value = 0.04978103906761788
$ws.Range("F8").Value = 1.059045595690989
$ws.Range("G8").Value = 1.017385126365411
$ws.Range("H8").Value = 0.7926862899448395
$ws.Range("I8").Value = 0.8040837796712168
$ws.Range("M8").Value = 1.175382807237639
$ws.Range("N8").Value = 1.579505822833738
$ws.Range("C9").Value = 0.000004364039200677183
$ws.Range("D9").Value = 0.01231531663410124
$ws.Range("E9").Value = 0.04978629060947171
$ws.Range("F9").Value = 1.268571307802631
$ws.Range("G9").Value = 1.262657887046629
$ws.Range("H9").Value = 0.8921979767808637
$ws.Range("I9").Value = 0.9647022589186065
$ws.Range("M9").Value = 1.5407802569257
$ws.Range("N9").Value = 1.899924622397236
$ws.Range("C10").Value = 0.000001001476743311969
$ws.Range("D10").Value = 0.01410156939377316
$ws.Range("E10").Value = 0.05007853355576941
$ws.Range("F10").Value = 1.427100737898428
$ws.Range("G10").Value = 1.447205716466556
$ws.Range("H10").Value = 0.9687475313120046
$ws.Range("I10").Value = 1.086219363514203
$ws.Range("M10").Value = 1.809146874145654
$ws.Range("N10").Value = 2.136173379601701
$ws.Range("C11").Value = 0.0000002824379321531012
$ws.Range("D11").Value = 0.01493241203350237
$ws.Range("E11").Value = 0.05027330041861688
$ws.Range("F11").Value = 1.500283625185205
$ws.Range("G11").Value = 1.532184486564859
$ws.Range("H11").Value = 1.004353823763267
$ws.Range("I11").Value = 1.142314936000446
$ws.Range("M11").Value = 1.931232462786454
$ws.Range("N11").Value = 2.243799872429918
$ws.Range("C12").Value = 0.0000001348037184101258
$ws.Range("D12").Value = 0.01524979541130023
$ws.Range("E12").Value = 0.05035587567177302
$ws.Range("F12").Value = 1.528154736086378
$ws.Range("G12").Value = 1.564517637188146
$ws.Range("H12").Value = 1.017952434103393
$ws.Range("I12").Value = 1.163678386517176
$ws.Range("M12").Value = 1.977464701674876
$ws.Range("N12").Value = 2.284574601459553
$ws.Range("C13").Value = 0.0000001609246100819917
$ws.Range("D13").Value = 0.01518131624654018
$ws.Range("E13").Value = 0.050337700369397
$ws.Range("F13").Value = 1.522145072561329
$ws.Range("G13").Value = 1.557547192032587
$ws.Range("H13").Value = 1.015018560838371
$ws.Range("I13").Value = 1.15907193002927
$ws.Range("M13").Value = 1.967507705822044
$ws.Range("N13").Value = 2.275792247643608
$ws.Range("C14").Value = 0.0000002677473402989961
$ws.Range("D14").Value = 0.01495846723521055
$ws.Range("E14").Value = 0.05027991753398098
$ws.Range("F14").Value = 1.502573396526572
$ws.Range("G14").Value = 1.534841443415644
$ws.Range("H14").Value = 1.005470264106521
$ws.Range("I14").Value = 1.14407006731895
$ws.Range("M14").Value = 1.93503599242473
$ws.Range("N14").Value = 2.247154074618493
$ws.Range("C15").Value = 0.0000003496553775406142
$ws.Range("D15").Value = 0.01482232941456374
$ws.Range("E15").Value = 0.05024567064143071
$ws.Range("F15").Value = 1.490605947180086
$ws.Range("G15").Value = 1.520953686752648
$ws.Range("H15").Value = 0.9996367478540265
$ws.Range("I15").Value = 1.1348969011338
$ws.Range("M15").Value = 1.915146285493734
$ws.Range("N15").Value = 2.229614727384728
$ws.Range("C16").Value = 0.00000106546426748011
$ws.Range("D16").Value = 0.01404765090300941
$ws.Range("E16").Value = 0.05006704321993638
$ws.Range("F16").Value = 1.422339952795738
$ws.Range("G16").Value = 1.441673318292004
$ws.Range("H16").Value = 0.9664365648615103
$ws.Range("I16").Value = 1.082570147319458
$ws.Range("M16").Value = 1.801168417535507
$ws.Range("N16").Value = 2.129142569809687
$ws.Range("C17").Value = 0.000001717893465258413
$ws.Range("D17").Value = 0.01357718440233668
$ws.Range("E17").Value = 0.04997324627325028
$ws.Range("F17").Value = 1.380737733549921
$ws.Range("G17").Value = 1.393304559359365
$ws.Range("H17").Value = 0.9462718483518699
$ws.Range("I17").Value = 1.050681282988421
$ws.Range("M17").Value = 1.731248086192039
$ws.Range("N17").Value = 2.067543545284309
$ws.Range("C18").Value = 0.000002168955577364073
$ws.Range("D18").Value = 0.01330829574266801
$ws.Range("E18").Value = 0.04992511997911109
$ws.Range("F18").Value = 1.356909375592465
$ws.Range("G18").Value = 1.365580442558468
$ws.Range("H18").Value = 0.9347472259857739
$ws.Range("I18").Value = 1.032416303682439
$ws.Range("M18").Value = 1.691032130464777
$ws.Range("N18").Value = 2.032128279290475
$ws.Range("C19").Value = 0.000002334453022889704
$ws.Range("D19").Value = 0.01321754450622592
$ws.Range("E19").Value = 0.04990982750630835
$ws.Range("F19").Value = 1.348858553957228
$ws.Range("G19").Value = 1.356209898400635
$ws.Range("H19").Value = 0.9308577470840476
$ws.Range("I19").Value = 1.026245145399528
$ws.Range("M19").Value = 1.67741574946254
$ws.Range("N19").Value = 2.020139938888121
$ws.Range("C20").Value = 0.000001640538906499955
$ws.Range("D20").Value = 0.01362708829238102
$ws.Range("E20").Value = 0.04998262893982286
$ws.Range("F20").Value = 1.385155961669284
$ws.Range("G20").Value = 1.39844348613056
$ws.Range("H20").Value = 0.948410781878323
$ws.Range("I20").Value = 1.054067946414307
$ws.Range("M20").Value = 1.738691183464198
$ws.Range("N20").Value = 2.074099353724989
$ws.Range("C21").Value = 0.0000002329218986574233
$ws.Range("D21").Value = 0.01502384736780016
$ws.Range("E21").Value = 0.05029665089014657
$ws.Range("F21").Value = 1.508317735445104
$ws.Range("G21").Value = 1.541506458202718
$ws.Range("H21").Value = 1.00827168058936
$ws.Range("I21").Value = 1.148473155955529
$ws.Range("M21").Value = 1.944573688519739
$ws.Range("N21").Value = 2.255565314398893
$ws.Range("C22").Value = 0.00000004342262860390633
$ws.Range("D22").Value = 0.01595286918922767
$ws.Range("E22").Value = 0.05055326835550034
$ws.Range("F22").Value = 1.589736334002254
$ws.Range("G22").Value = 1.63590408239736
$ws.Range("H22").Value = 1.048067425558884
$ws.Range("I22").Value = 1.210881117160866
$ws.Range("M22").Value = 2.07913634713384
$ws.Range("N22").Value = 2.374272474789848
$ws.Range("C23").Value = 0.00000007499756771167654
$ws.Range("D23").Value = 0.01545550918010719
$ws.Range("E23").Value = 0.05041162743387062
$ws.Range("F23").Value = 1.546195393572447
$ws.Range("G23").Value = 1.585438196227472
$ws.Range("H23").Value = 1.02676521544339
$ws.Range("I23").Value = 1.177506689885476
$ws.Range("M23").Value = 2.007316953280679
$ws.Range("N23").Value = 2.310907419963826
$ws.Range("C24").Value = 0.000001675274852619424
$ws.Range("D24").Value = 0.01360452180549032
$ws.Range("E24").Value = 0.04997836896900409
$ws.Range("F24").Value = 1.383158202918764
$ws.Range("G24").Value = 1.396119916987288
$ws.Range("H24").Value = 0.947443557358838
$ws.Range("I24").Value = 1.052536623260991
$ws.Range("M24").Value = 1.735326214744248
$ws.Range("N24").Value = 2.071135476096288
$ws.Range("C25").Value = 0.000006268529575859105
$ws.Range("D25").Value = 0.01167282175692463
$ws.Range("E25").Value = 0.04973394683104004
$ws.Range("F25").Value = 1.211104831029076
$ws.Range("G25").Value = 1.195568090093758
$ws.Range("H25").Value = 0.8646852273309662
$ws.Range("I25").Value = 0.9206512817370935
$ws.Range("M25").Value = 1.441951971154381
$ws.Range("N25").Value = 1.813084641365805
